$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 13:05"

# Finlandia (row 58): new cases pushed totals/active up
$ws.Cells.Item(58, 2).Value = 6054
$ws.Cells.Item(58, 3).Value = 51
$ws.Cells.Item(58, 5).Value = 1479

# Bosnia y Herzegovina (row 78)
$ws.Cells.Item(78, 2).Value = 2181
$ws.Cells.Item(78, 3).Value = 23
$ws.Cells.Item(78, 4).Value = 1228
$ws.Cells.Item(78, 5).Value = 833
$ws.Cells.Item(78, 7).Value = 3
$ws.Cells.Item(78, 8).Value = 120

# Eslovenia (row 91)
$ws.Cells.Item(91, 2).Value = 1463
$ws.Cells.Item(91, 3).Value = 2
$ws.Cells.Item(91, 4).Value = 260
$ws.Cells.Item(91, 7).Value = 1
$ws.Cells.Item(91, 8).Value = 103

# Libano (row 106)
$ws.Cells.Item(106, 2).Value = 878
$ws.Cells.Item(106, 3).Value = 8
$ws.Cells.Item(106, 4).Value = 236
$ws.Cells.Item(106, 5).Value = 616

# Malta gained cases and now overtakes Jamaica in the ranking, so the two
# rows swap places (row 123 becomes Malta, row 124 becomes Jamaica).
$ws.Cells.Item(123, 1).Value = "Malta"
$ws.Cells.Item(123, 2).Value = 508
$ws.Cells.Item(123, 3).Value = 2
$ws.Cells.Item(123, 4).Value = 436
$ws.Cells.Item(123, 5).Value = 66
$ws.Cells.Item(123, 6).Value = 1
$ws.Cells.Item(123, 7).Value = 1
$ws.Cells.Item(123, 8).Value = 6

$ws.Cells.Item(124, 1).Value = "Jamaica"
$ws.Cells.Item(124, 2).Value = 507
$ws.Cells.Item(124, 3).Value = 2
$ws.Cells.Item(124, 4).Value = 100
$ws.Cells.Item(124, 5).Value = 398
$ws.Cells.Item(124, 6).Value = 0
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 9

# Nepal (row 142)
$ws.Cells.Item(142, 4).Value = 35
$ws.Cells.Item(142, 5).Value = 184
